# Update cryptos list (GitHub Actions scheduled refresh).
# Updates Price (column D) and Volume(1h) (column E) for each ranked coin,
# and also reflects the rank swap between Cosmos and OKB (rows 33/34).
# Numeric-looking Price strings are apostrophe-prefixed so Excel keeps them
# as text (matching the original inline-string cell type) instead of
# auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.608.92"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "3.572.00"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'587.95"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'182.52"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "3.566.81"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "'0.996"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'0.670"
$ws.Range("E10").Value = "  -6.16%  "
$ws.Range("D11").Value = "'53.85"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "'0.144"
$ws.Range("E12").Value = "  -10.44%  "
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  -13.78%  "
$ws.Range("D14").Value = "'9.77"
$ws.Range("E14").Value = "  -7.93%  "
$ws.Range("D15").Value = "4.134.24"
$ws.Range("E15").Value = "  -3.12%  "
$ws.Range("D16").Value = "3.565.02"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'18.31"
$ws.Range("E18").Value = "  -5.30%  "
$ws.Range("D19").Value = "66.338.80"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").Value = "'12.09"
$ws.Range("E20").Value = "  -5.67%  "
$ws.Range("E21").Value = "  -6.40%  "
$ws.Range("D22").Value = "'393.01"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("D23").Value = "'4.28"
$ws.Range("E23").Value = "  -6.50%  "
$ws.Range("D24").Value = "'84.53"
$ws.Range("E24").Value = "  -4.61%  "
$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = "  -5.10%  "
$ws.Range("D26").Value = "'12.32"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").Value = "'6.02"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "'10.23"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").Value = "'3.58"
$ws.Range("E29").Value = "  -12.17%  "
$ws.Range("D30").Value = "'8.91"
$ws.Range("E30").Value = "  -7.12%  "
$ws.Range("D31").Value = "'31.05"
$ws.Range("E31").Value = "  -5.47%  "
$ws.Range("D32").Value = "'6.75"
$ws.Range("E32").Value = "  -7.92%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'65.47"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.92"
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("D35").Value = "'603.86"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  -6.21%  "
$ws.Range("D37").Value = "'41.33"
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'0.372"
$ws.Range("E40").Value = "  -7.09%  "
$ws.Range("D41").Value = "0.0₃0743"
$ws.Range("E41").Value = "  -14.79%  "
$ws.Range("E42").Value = "  -6.11%  "
$ws.Range("D43").Value = "2.912.69"
$ws.Range("E43").Value = "  +6.85%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("E44").Value = "  -8.57%  "
$ws.Range("D45").Value = "'0.0406"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  -8.86%  "
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").Value = "'3.04"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "'136.23"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").Value = "'2.51"
$ws.Range("E50").Value = "  -8.37%  "
$ws.Range("D51").Value = "'8.23"
$ws.Range("E51").Value = "  -8.69%  "
